# Update recomputed NATMI TPM-derived metrics for Dhh-Ptch2 LR-pair sheet.
# New TPM normalization changed ligand/receptor expressing-cell counts and all
# downstream average/total expression, specificity, and edge-weight columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9167996666666666
$ws.Range("H2").Value = 2.750399
$ws.Range("I2").Value = 0.2833456974325495
$ws.Range("J2").Value = 0.2833456974325495
$ws.Range("M2").Value = 1.418869666666667
$ws.Range("N2").Value = 4.256609
$ws.Range("O2").Value = 0.1110803932403577
$ws.Range("P2").Value = 0.1110803932403577
$ws.Range("Q2").Value = 1.300819237443444
$ws.Range("R2").Value = 11.707373136991
$ws.Range("S2").Value = 0.03147415149377101
$ws.Range("T2").Value = 0.03147415149377101

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9167996666666666
$ws.Range("H3").Value = 2.750399
$ws.Range("I3").Value = 0.2833456974325495
$ws.Range("J3").Value = 0.2833456974325495
$ws.Range("O3").Value = 0.7197921976927236
$ws.Range("P3").Value = 0.7197921976927237
$ws.Range("Q3").Value = 8.429206184878776
$ws.Range("R3").Value = 75.862855663909
$ws.Range("S3").Value = 0.2039500222617523
$ws.Range("T3").Value = 0.2039500222617524

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9167996666666666
$ws.Range("H4").Value = 2.750399
$ws.Range("I4").Value = 0.2833456974325495
$ws.Range("J4").Value = 0.2833456974325495
$ws.Range("N4").Value = 6.480974999999999
$ws.Range("O4").Value = 0.1691274090669186
$ws.Range("P4").Value = 0.1691274090669185
$ws.Range("Q4").Value = 1.980585239891666
$ws.Range("R4").Value = 17.825267159025
$ws.Range("S4").Value = 0.04792152367702614
$ws.Range("T4").Value = 0.04792152367702613

# Row 5
$ws.Range("I5").Value = 0.2271242616180895
$ws.Range("J5").Value = 0.2271242616180895
$ws.Range("M5").Value = 1.418869666666667
$ws.Range("N5").Value = 4.256609
$ws.Range("O5").Value = 0.1110803932403577
$ws.Range("P5").Value = 0.1110803932403577
$ws.Range("Q5").Value = 1.042710764553889
$ws.Range("R5").Value = 9.384396880984999
$ws.Range("S5").Value = 0.02522905229496326
$ws.Range("T5").Value = 0.02522905229496326

# Row 6
$ws.Range("I6").Value = 0.2271242616180895
$ws.Range("J6").Value = 0.2271242616180895
$ws.Range("O6").Value = 0.7197921976927236
$ws.Range("P6").Value = 0.7197921976927237
$ws.Range("R6").Value = 60.810152520515
$ws.Range("S6").Value = 0.1634822714194218
$ws.Range("T6").Value = 0.1634822714194218

# Row 7
$ws.Range("I7").Value = 0.2271242616180895
$ws.Range("J7").Value = 0.2271242616180895
$ws.Range("N7").Value = 6.480974999999999
$ws.Range("O7").Value = 0.1691274090669186
$ws.Range("P7").Value = 0.1691274090669185
$ws.Range("S7").Value = 0.03841293790370445
$ws.Range("T7").Value = 0.03841293790370445

# Row 8
$ws.Range("G8").Value = 1.583934333333334
$ws.Range("H8").Value = 4.751803000000001
$ws.Range("I8").Value = 0.4895300409493609
$ws.Range("J8").Value = 0.4895300409493609
$ws.Range("M8").Value = 1.418869666666667
$ws.Range("N8").Value = 4.256609
$ws.Range("O8").Value = 0.1110803932403577
$ws.Range("P8").Value = 0.1110803932403577
$ws.Range("Q8").Value = 2.247396379558556
$ws.Range("R8").Value = 20.226567416027
$ws.Range("S8").Value = 0.05437718945162341
$ws.Range("T8").Value = 0.05437718945162341

# Row 9
$ws.Range("G9").Value = 1.583934333333334
$ws.Range("H9").Value = 4.751803000000001
$ws.Range("I9").Value = 0.4895300409493609
$ws.Range("J9").Value = 0.4895300409493609
$ws.Range("O9").Value = 0.7197921976927236
$ws.Range("P9").Value = 0.7197921976927237
$ws.Range("S9").Value = 0.3523599040115495
$ws.Range("T9").Value = 0.3523599040115495

# Row 10
$ws.Range("G10").Value = 1.583934333333334
$ws.Range("H10").Value = 4.751803000000001
$ws.Range("I10").Value = 0.4895300409493609
$ws.Range("J10").Value = 0.4895300409493609
$ws.Range("N10").Value = 6.480974999999999
$ws.Range("O10").Value = 0.1691274090669186
$ws.Range("P10").Value = 0.1691274090669185
$ws.Range("S10").Value = 0.08279294748618796
$ws.Range("T10").Value = 0.08279294748618794
